# Insert a new data row at row 238 (pushing existing rows 238:249 down to
# 239:250) and populate it with a new "Perejil" price observation for
# Feria Lagunitas de Puerto Montt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 238; this shifts rows 238-249
# down to 239-250 (and extends the sheet dimension to A1:R250).
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new record.
$ws.Cells.Item(238, 1).Value = 4
$ws.Cells.Item(238, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(238, 3).Value = 'Los Lagos'
$ws.Cells.Item(238, 4).Value = 44706
$ws.Cells.Item(238, 5).Value = 10
$ws.Cells.Item(238, 6).Value = 100112044
$ws.Cells.Item(238, 7).Value = 'Perejil'
$ws.Cells.Item(238, 8).Value = 'Sin especificar'
$ws.Cells.Item(238, 9).Value = 'Primera'
$ws.Cells.Item(238, 10).Value = 25
$ws.Cells.Item(238, 11).Value = 6000
$ws.Cells.Item(238, 12).Value = 6000
$ws.Cells.Item(238, 13).Value = 6000
$ws.Cells.Item(238, 14).Value = '$/docena de atados (2 kilos)'
$ws.Cells.Item(238, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(238, 16).Value = 3000
$ws.Cells.Item(238, 17).Value = 2
$ws.Cells.Item(238, 18).Value = 'Hortaliza'

# Keep the date column's number format consistent with the rest of column D.
$ws.Cells.Item(238, 4).NumberFormat = $ws.Cells.Item(239, 4).NumberFormat
